# Apply schedule regeneration changes to the "Horario" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing time-block labels in column A (rows 2-8) ---
$ws.Range("A2").Value = "08:31-09:50"
$ws.Range("A3").Value = "10:01-10:40"
$ws.Range("A4").Value = "10:41-11:20"
$ws.Range("A5").Value = "11:31-12:10"
$ws.Range("A6").Value = "13:41-14:20"
$ws.Range("A7").Value = "14:31-15:10"
$ws.Range("A8").Value = "16:15-17:45"

# --- Fill in newly scheduled course sections on row 3 ---
$ws.Range("B3").Value = "DSY1105-004D"
$ws.Range("C3").Value = "DSY1104-003D"
$ws.Range("E3").Value = "DSY1104-003D"

# --- Append four new time blocks (rows 9-12) ---
$newBlocks = @(
    @{ Row = 9;  Time = "18:00-19:30" },
    @{ Row = 10; Time = "19:01-20:20" },
    @{ Row = 11; Time = "20:31-21:10" },
    @{ Row = 12; Time = "21:11-22:30" }
)

foreach ($block in $newBlocks) {
    $r = $block.Row
    $ws.Range("A$r").Value = $block.Time
    # Copy the (empty) day-cell formatting from row 2 so the new row gets the
    # same style (vertical=top, wrap text) as the rest of the schedule grid.
    $ws.Range("B2:G2").Copy($ws.Range("B$r`:G$r"))
}
